# CalculatorAdd.feature.xlsx update:
# Replace the "SummandOne"/"SummandTwo" placeholders used in the Gherkin
# Examples table with the friendlier "Number1"/"Number2" names (and their
# <Number1>/<Number2> placeholder forms) on both example sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Add two positive numbers")
$ws1.Range("C1").Value = "<Number1>"
$ws1.Range("C2").Value = "<Number2>"
$ws1.Range("B7").Value = "Number1"
$ws1.Range("C7").Value = "Number2"

$ws2 = $wb.Worksheets.Item("Add two negative numbers")
$ws2.Range("C1").Value = "<Number1>"
$ws2.Range("C2").Value = "<Number2>"
$ws2.Range("B7").Value = "Number1"
$ws2.Range("C7").Value = "Number2"

# Restore the cell selections recorded in each sheet view, and make the
# first sheet ("Add two positive numbers") the active tab.
$ws2.Range("C8").Select()
$ws1.Activate()
$ws1.Range("C3").Select()
